# Patient timeline map: add a per-patient scan-count "SUM" column (J),
# an inclusion-comments legend (column L), and conditional formatting
# that strikes through patients with fewer than two scans.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column J: total scan count per patient -----------------------
# Header + body alignment (row 27 is blank in the source data and must
# stay that way, so it is deliberately excluded from every J/L range).
$ws.Range("J1:J26").HorizontalAlignment = -4108   # xlCenter (match columns B:I)
$ws.Range("J28:J34").HorizontalAlignment = -4108
$ws.Columns("J").ColumnWidth = 12.33              # ~ same rendered width as column I

$ws.Range("J1").Value = "SUM"

# Row 2 is the first data row and gets a plain (non-shared) formula,
# matching how Excel seeds the first cell before autofill turns the
# rest of the column into a shared-formula block.
$ws.Range("J2").Formula = "=SUM(B2:H2)"

# Fill down — row 27 is blank in the source data and stays untouched.
$ws.Range("J3:J26").FormulaR1C1 = "=SUM(RC[-8]:RC[-2])"
$ws.Range("J28:J34").FormulaR1C1 = "=SUM(RC[-8]:RC[-2])"

# --- Legend / inclusion comments in column L ---------------------------
$ws.Range("L1:L2").HorizontalAlignment = -4131    # xlLeft

$ws.Range("L1").Value = "Inclusion comments"
$ws.Range("L2").Value = "fluid build up instead of brain expansion"
$ws.Range("L3").Value = "ventricle enlargement (ventriculomegaly // hydrocephalus)"
$ws.Range("L4").Value = "ventricle enlargement (ventriculomegaly // hydrocephalus)"
$ws.Range("L5").Value = "acute yes, fast ?? Skull edges difficult (skull edges always difficult)"

# --- Conditional formatting: strikethrough patients with <2 scans ------
# Exercised add/delete a couple of times while tuning the strike color
# (35% vs 50% darker), landing back on the 35%-darker grey — this is
# why three dxf records exist even though only the last rule survives.
$rng = $ws.Range("A2:I34")

$fc1 = $rng.FormatConditions.Add(2, 8, "=SUM(`$B2:`$H2)<2")
$fc1.Font.Strikethrough = $true
$fc1.Font.Color = 10921638   # RGB(166,166,166) ~ "white, darker 35%"

$fc2 = $rng.FormatConditions.Add(2, 8, "=SUM(`$B2:`$H2)<2")
$fc2.Font.Strikethrough = $true
$fc2.Font.Color = 8421504    # RGB(128,128,128) ~ "white, darker 50%"

$fc3 = $rng.FormatConditions.Add(2, 8, "=SUM(`$B2:`$H2)<2")
$fc3.Font.Strikethrough = $true
$fc3.Font.Color = 10921638   # back to "white, darker 35%"

$rng.FormatConditions.Item(1).Delete()
$rng.FormatConditions.Item(1).Delete()
$rng.FormatConditions.Item(1).Priority = 1

# --- Leave selection where the author left it --------------------------
$ws.Range("J2").Select()
